$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'89.008.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "'3.163.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.38%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'210.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").Value = "'612.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.17%  "

$ws.Range("D7").Value = "'0.381"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.82%  "

$ws.Range("D8").Value = "'0.683"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.55%  "

$ws.Range("D9").Value = "'0.998"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("D10").Value = "'3.149.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.57%  "

$ws.Range("D11").Value = "'0.570"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.68%  "

$ws.Range("E12").Value = "  -5.82%  "

$ws.Range("D13").Value = "'0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.82%  "

$ws.Range("D14").Value = "'88.677.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "

$ws.Range("D15").Value = "'3.725.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.83%  "

$ws.Range("D16").Value = "'5.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.06%  "

$ws.Range("D17").Value = "'32.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.15%  "

$ws.Range("D18").Value = "'3.134.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.90%  "

$ws.Range("D19").Value = "'3.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.74%  "

$ws.Range("D20").Value = "'13.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.54%  "

$ws.Range("D21").Value = "'432.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "

$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").Value = "'0.0000185"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +37.01%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'8.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.37%  "

$ws.Range("D24").Value = "'5.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.61%  "

$ws.Range("D25").Value = "'5.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.40%  "

$ws.Range("D26").Value = "'11.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.43%  "

$ws.Range("D27").Value = "'3.353.69"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").Value = "'74.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.45%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").Value = "'0.167"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.33%  "

$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").Value = "'3.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +29.17%  "

$ws.Range("D33").Value = "'8.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.44%  "

$ws.Range("D34").Value = "'525.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.55%  "

$ws.Range("D35").Value = "'6.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.91%  "

$ws.Range("D36").Value = "'1.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.20%  "

$ws.Range("D37").Value = "'1.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.86%  "

$ws.Range("D38").Value = "'22.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.81%  "

$ws.Range("D39").Value = "'21.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.47%  "

$ws.Range("D40").Value = "'0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.21%  "

$ws.Range("D41").Value = "'0.126"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.93%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").Value = "'1.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.81%  "

$ws.Range("D44").Value = "'0.367"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.54%  "

$ws.Range("D45").Value = "'150.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.25%  "

$ws.Range("D46").Value = "'43.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.16%  "

$ws.Range("D47").Value = "'169.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.02%  "

$ws.Range("E48").Value = "  -8.16%  "

$ws.Range("D49").Value = "'1.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.13%  "

$ws.Range("D50").Value = "'0.603"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.56%  "

$ws.Range("D51").Value = "'4.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.84%  "
